# Refactored Dispatcher to use List<MailMessage> as TransactionItem and
# removed Orchestrator transaction dependency.
#
# Functionally this appends two new "Partnership Offer" rows (25 and 26) to
# the report sheet, reproducing an existing row's Sender/Company/Address/
# Email/Subject/VAT-ID values and formatting, with new "Date Processed"
# timestamps.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Partnership_Emails")

$lastRow = 24

# Duplicate the last existing data row (including its number formatting)
# into the two new rows so that formatting/style indexes line up with the
# rest of the sheet.
$ws.Range("A" + $lastRow + ":G" + $lastRow).Copy($ws.Range("A25:G25"))
$ws.Range("A" + $lastRow + ":G" + $lastRow).Copy($ws.Range("A26:G26"))

# Update the "Date Processed" timestamps for the two new rows.
$ws.Range("A25").Value2 = 45857.012384259258
$ws.Range("A26").Value2 = 45857.758298611108
